# Apply the "cryptos list" refresh (GitHub Actions scheduled update).
# Updates Price (D) / Volume(1h) (E) figures for most rows, and for three
# row-pairs the ranking flipped so Coin (B) / Link (C) / Price (D) / Volume (E)
# were all swapped between the two adjacent rows (15<->16, 38<->39, 43<->44).
#
# All of these columns are stored as text in the sheet (prices use "."
# as a thousands separator in many rows, e.g. "55.796.99", so they can
# never be read back as numbers). For the handful of price values that
# *would* parse as a plain number (e.g. "504.16"), a leading apostrophe
# is used so Excel keeps them as literal text instead of silently
# converting the cell to a Number, matching every other cell in the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '55.796.99'
$ws.Range("E2").Value = '  -1.62%  '

$ws.Range("D3").Value = '2.345.16'
$ws.Range("E3").Value = '  -1.89%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = "'504.16"
$ws.Range("E5").Value = '  -0.12%  '

$ws.Range("D6").Value = "'129.03"
$ws.Range("E6").Value = '  -2.56%  '

$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("E8").Value = '  -2.44%  '

$ws.Range("D9").Value = '2.355.96'
$ws.Range("E9").Value = '  -1.52%  '

$ws.Range("D10").Value = "'0.0972"
$ws.Range("E10").Value = '  -0.27%  '

$ws.Range("E11").Value = '  +0.00%  '

$ws.Range("D12").Value = "'4.78"
$ws.Range("E12").Value = '  +2.86%  '

$ws.Range("E13").Value = '  -1.25%  '

$ws.Range("D14").Value = '2.761.73'
$ws.Range("E14").Value = '  -1.91%  '

# Row 15/16: ranking swap (WrappedBTC now ranks above Avalanche)
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '55.723.93'
$ws.Range("E15").Value = '  -1.51%  '

$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").Value = "'21.62"
$ws.Range("E16").Value = '  -0.06%  '

$ws.Range("E17").Value = '  -0.76%  '

$ws.Range("D18").Value = '2.315.79'
$ws.Range("E18").Value = '  -3.83%  '

$ws.Range("E19").Value = '  -2.85%  '

$ws.Range("D20").Value = "'311.01"
$ws.Range("E20").Value = '  +0.67%  '

$ws.Range("E21").Value = '  -1.66%  '

$ws.Range("E22").Value = '  -0.54%  '

$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = '  -0.23%  '

$ws.Range("D24").Value = "'65.29"
$ws.Range("E24").Value = '  -2.65%  '

$ws.Range("D25").Value = "'0.997"
$ws.Range("E25").Value = '  -0.40%  '

$ws.Range("E26").Value = '  -1.48%  '

$ws.Range("E27").Value = '  -2.73%  '

$ws.Range("D28").Value = "'7.09"
$ws.Range("E28").Value = '  -4.28%  '

$ws.Range("D29").Value = "'171.49"
$ws.Range("E29").Value = '  -2.36%  '

$ws.Range("D30").Value = "'1.64"
$ws.Range("E30").Value = '  -0.88%  '

$ws.Range("D31").Value = '0.0₃0703'
$ws.Range("E31").Value = '  -2.94%  '

$ws.Range("E32").Value = '  -0.03%  '

$ws.Range("E33").Value = '  -1.75%  '

$ws.Range("E34").Value = '  +0.03%  '

$ws.Range("D35").Value = "'1.07"
$ws.Range("E35").Value = '  -5.33%  '

$ws.Range("D36").Value = "'17.69"
$ws.Range("E36").Value = '  -0.87%  '

$ws.Range("E37").Value = '  -1.97%  '

# Row 38/39: ranking swap (SuiNetwork now ranks above NEARProtocol)
$ws.Range("B38").Value = 'SuiNetwork'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D38").Value = "'0.829"
$ws.Range("E38").Value = '  +0.20%  '

$ws.Range("B39").Value = 'NEARProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D39").Value = "'3.63"
$ws.Range("E39").Value = '  -4.72%  '

$ws.Range("D40").Value = "'36.01"
$ws.Range("E40").Value = '  -2.14%  '

$ws.Range("E41").Value = '  -4.00%  '

$ws.Range("E42").Value = '  -1.09%  '

# Row 43/44: ranking swap (RenderToken now ranks above Aave)
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = "'4.88"
$ws.Range("E43").Value = '  +0.67%  '

$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").Value = "'126.48"
$ws.Range("E44").Value = '  -3.87%  '

$ws.Range("D45").Value = "'0.555"
$ws.Range("E45").Value = '  -1.99%  '

$ws.Range("E46").Value = '  -2.01%  '

$ws.Range("D47").Value = "'238.52"
$ws.Range("E47").Value = '  -4.97%  '

$ws.Range("E48").Value = '  -1.94%  '

$ws.Range("D50").Value = "'16.68"
$ws.Range("E50").Value = '  -2.50%  '

$ws.Range("D51").Value = "'0.952"
$ws.Range("E51").Value = '  +0.06%  '
